$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B,C,D,E,F,H,I,J,K,L,M,N for rows 2-25 (column G and A/O are unchanged)
$colIndexes = @(2, 3, 4, 5, 6, 8, 9, 10, 11, 12, 13, 14)

$data = @(
    @(1.41439309687118, 0.07513119594497653, 0.0116182754094396, 0.06358328784659673, 3.767704608165189, 0.07973214163530429, 2.378409290992209, 0.1458084058800759, 1.245497301426269, 0.3418581033704982, 0.3636948708781311, 4.004334293927499),
    @(1.385652285262694, 0.06958702018067697, 0.01145343688971856, 0.06374172123274935, 3.769346261606387, 0.07973214163530429, 2.382981102763537, 0.1463132986270743, 1.212520268991852, 0.3407908606219223, 0.3590000546155778, 4.024153980476456),
    @(1.36873286538119, 0.06622488734093679, 0.01135633486415877, 0.06385006569862028, 3.771710713770716, 0.07973214163530429, 2.386574929847782, 0.1466427731050635, 1.192950057888055, 0.3402659181184831, 0.3562897929574689, 4.037206725149879),
    @(1.362021437493667, 0.06486530894530063, 0.01131780850657904, 0.06389700823586697, 3.773015631818694, 0.07973214163530429, 2.388237422814512, 0.1467819426260775, 1.185145767169928, 0.3400848518814286, 0.3552287840750807, 4.042748035554368),
    @(1.360918097843665, 0.06464018590193632, 0.01131147455966364, 0.06390497185234167, 3.773252939854117, 0.07973214163530429, 2.388525440517625, 0.1468053482395923, 1.183860190781559, 0.3400567723660615, 0.3550552312675705, 4.043681589290244),
    @(1.368641609684005, 0.06620650906445746, 0.01135581104564309, 0.06385068746847633, 3.771726929723116, 0.07973214163530429, 2.386596549001737, 0.1466446301160014, 1.192844114838067, 0.3402633430793429, 0.356275307781825, 4.037280557383774),
    @(1.404332501681523, 0.07321081431493326, 0.01156059212556215, 0.06363562419010105, 3.7679892804031, 0.07973214163530429, 2.379822457075463, 0.1459784586039472, 1.233986296130496, 0.341463107890803, 0.3620403792188398, 4.010984824322193),
    @(1.480081863559917, 0.08728246959114472, 0.0119943642857443, 0.0633013118491661, 3.771413321634924, 0.07973214163530429, 2.372775455262939, 0.1448260928000966, 1.320037491808108, 0.3448472720327089, 0.3747100147999305, 3.966424602519595),
    @(1.539236436787519, 0.09783098796006584, 0.01233218623664811, 0.06310849629422055, 3.780476821145058, 0.07973214163530429, 2.371395330883999, 0.1440726647764965, 1.38653359748227, 0.3479592734737764, 0.3848471188867961, 3.937952070596396),
    @(1.566906207618956, 0.1026765632288118, 0.01248992391104053, 0.06303214242175059, 3.786020607849537, 0.07973214163530429, 2.371591185374257, 0.1437500110240304, 1.417495914043172, 0.3495102468378519, 0.3896381665626478, 3.925924507755965),
    @(1.577493034188478, 0.1045182850908475, 0.01255022999182742, 0.06300485413346646, 3.788324018584575, 0.07973214163530429, 2.371783695792601, 0.1436307079409715, 1.429322913603443, 0.3501169541936662, 0.3914781673294812, 3.921502910547503),
    @(1.575208135486321, 0.1041213337798865, 0.0125372166111255, 0.06301065899919056, 3.787818863707599, 0.07973214163530429, 2.371736973379143, 0.1436562741012679, 1.426771214694753, 0.3499854277547598, 0.3910807471823077, 3.922449266193837),
    @(1.567775012338814, 0.1028279460768715, 0.01249487387514137, 0.06302986487150353, 3.786206021541346, 0.07973214163530429, 2.371604652082382, 0.1437401382429879, 1.418466880879691, 0.3495597728623494, 0.3897890293091564, 3.925558075141851),
    @(1.563236177670888, 0.1020365967800103, 0.01246901222553731, 0.06304184044029792, 3.785244682885235, 0.07973214163530429, 2.371539010511782, 0.1437918820957851, 1.413393543551081, 0.3493015697003869, 0.389001163452356, 3.927479628473321),
    @(1.53744341134302, 0.09751526818632783, 0.01232195859034491, 0.06311371407302069, 3.78014309334867, 0.07973214163530429, 2.371399097304497, 0.1440941543176262, 1.38452445937935, 0.3478606304910983, 0.384537619067757, 3.938756708394322),
    @(1.521814760772827, 0.09475365245363321, 0.01223277906379039, 0.06316070969878851, 3.777377176514491, 0.07973214163530429, 2.371524157595339, 0.1442847260495732, 1.366996608479639, 0.3470112668621397, 0.3818453182075316, 3.945911680021268),
    @(1.51289713901167, 0.09316966517735636, 0.01218186814768174, 0.06318880997483056, 3.775920034399462, 0.07973214163530429, 2.371673623390151, 0.1443962288579037, 1.356982170309152, 0.3465354753658332, 0.3803136910710023, 3.950114070315621),
    @(1.509890089377024, 0.09263411213514416, 0.01216469662792008, 0.06319850820475992, 3.775449647402525, 0.07973214163530429, 2.371737549639185, 0.1444343068544409, 1.353602989561381, 0.3463765707809756, 0.3797980161985564, 3.951551876601485),
    @(1.523471053613406, 0.09504717324749379, 0.01224223283046655, 0.06315559629019774, 3.77765777296608, 0.07973214163530429, 2.37150282091477, 0.1442642437203929, 1.368855533644279, 0.3471003650357005, 0.3821301684794918, 3.945141012748309),
    @(1.569955348942358, 0.1032076602482732, 0.01250729545467166, 0.06302417959648121, 3.786674214794928, 0.07973214163530429, 2.371640306997278, 0.1437154272525536, 1.420903292354183, 0.3496842724076288, 0.3901677404762225, 3.924641333507424),
    @(1.600969999031122, 0.1085807134528523, 0.01268387197209009, 0.06294776075781705, 3.793756534590699, 0.07973214163530429, 2.372419920501343, 0.1433735201173132, 1.455515332997493, 0.3514859778497907, 0.3955707178914949, 3.912018721001672),
    @(1.584358976894748, 0.1057093647558531, 0.01258932705634308, 0.06298768313286907, 3.789867789205232, 0.07973214163530429, 2.371940745398177, 0.1435544703980547, 1.436987807138621, 0.3505140596673328, 0.3926733571748855, 3.918684721568269),
    @(1.522722033337999, 0.09491446096841116, 0.01223795766197, 0.06315790469205318, 3.777530500928634, 0.07973214163530429, 2.371512225592234, 0.1442734977359255, 1.368014918697895, 0.3470600447511885, 0.3820013371638424, 3.945489154474203),
    @(1.458974061452722, 0.08343919201830374, 0.01187362266984238, 0.06338244389146563, 3.769336608620165, 0.07973214163530429, 2.374014663730293, 0.1451214210623206, 1.296183350907938, 0.3438216699378174, 0.3711367947435278, 3.977729850542161)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $rowNum = $startRow + $i
    for ($j = 0; $j -lt $colIndexes.Length; $j++) {
        $ws.Cells.Item($rowNum, $colIndexes[$j]).Value = $row[$j]
    }
}
